# Update the PEBCOM interactive map workbook.
# The row for "Vuelta de Obligado 2775" (row 40) was removed upstream,
# so all subsequent rows shift up by one and the sheet shrinks by a row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PEBCOM")

# Delete entire row 40; Excel automatically shifts rows 41:82 up to 40:81.
$ws.Rows.Item(40).Delete()

$wb.Save()
